$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 34; this shifts the existing rows 34-84 down to 35-85,
# carrying all of their original values/styles with them.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with a fresh weekly record (same market /
# product / pricing details as the record that used to sit at row 34, just a
# week later in date).
$ws.Range("A34").Value = 5
$ws.Range("B34").Value = "Macroferia Regional de Talca"
$ws.Range("C34").Value = "Maule"
$ws.Range("D34").Value = 44883
$ws.Range("E34").Value = 7
$ws.Range("F34").Value = 100112040
$ws.Range("G34").Value = "Cilantro"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 150
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = 7000
$ws.Range("N34").Value = "$/caja 36 atados"
$ws.Range("O34").Value = "Región del Maule"
$ws.Range("P34").Value = 194
$ws.Range("Q34").Value = 36
$ws.Range("R34").Value = "Hortaliza"
